$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team")

# --- Column B width (Player names) ---
$ws.Columns("B:B").ColumnWidth = 16.67

# --- New "Injured" column (L) ---
$ws.Range("L1").Value = "Injured"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").HorizontalAlignment = -4108
$ws.Range("L1").VerticalAlignment = -4108
$ws.Range("L1").WrapText = $true

$ws.Range("L2").Value = $false
$ws.Range("L3").Value = $false
$ws.Range("L4").Value = $false
$ws.Range("L5").Value = $true
$ws.Range("L6").Value = $false
$ws.Range("L7").Value = $false
$ws.Range("L8").Value = $false
$ws.Range("L9").Value = $false
$ws.Range("L10").Value = $false
$ws.Range("L11").Value = $false
$ws.Range("L12").Value = $false
$ws.Range("L13").Value = $false
$ws.Range("L14").Value = $false
$ws.Range("L15").Value = $false

# --- A1 header style (id column - centered, bold, no wrap) ---
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").WrapText = $false
$excel.CutCopyMode = $false

# --- Trade logic data updates ---
# Free Transfers count bumped
$ws.Range("K2").Value = 2
# Jordan Pickford (row 12) no longer starting
$ws.Range("I12").Value = $false

# --- Row heights re-fit to the now-narrower/wider name column ---
# Rows whose wrapped text now only needs one line shrink back to the
# sheet's default row height ...
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(10).AutoFit()
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).AutoFit()
# ... while the long player names still wrap onto two lines
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30
